$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Place 1 row (row 3): 100 Meter Dash / 200 Meter Dash results ---
$t.Cell(3,2).Range.Text  = "P1"
$t.Cell(3,3).Range.Text  = "PLM"
$t.Cell(3,4).Range.Text  = "11.3"
$t.Cell(3,5).Range.Text  = "5"
$t.Cell(3,7).Range.Text  = "1"
$t.Cell(3,8).Range.Text  = "P5"
$t.Cell(3,9).Range.Text  = "PLM"
$t.Cell(3,10).Range.Text = "11.2"
$t.Cell(3,11).Range.Text = "0"
$t.Cell(3,12).Range.Text = "5"

# --- Place 2 row (row 4) ---
$t.Cell(4,2).Range.Text  = "P2"
$t.Cell(4,3).Range.Text  = "PLM"
$t.Cell(4,4).Range.Text  = "11.4"
$t.Cell(4,5).Range.Text  = "3"
$t.Cell(4,7).Range.Text  = "2"
$t.Cell(4,8).Range.Text  = "P2"
$t.Cell(4,9).Range.Text  = "PLM"
$t.Cell(4,10).Range.Text = "11.4"
$t.Cell(4,11).Range.Text = "3"
$t.Cell(4,12).Range.Text = "0"

# --- Place 3 row (row 5) ---
$t.Cell(5,2).Range.Text  = "G1"
$t.Cell(5,3).Range.Text  = "GWY"
$t.Cell(5,4).Range.Text  = "11.5"
$t.Cell(5,6).Range.Text  = "1"
$t.Cell(5,7).Range.Text  = "3"
$t.Cell(5,8).Range.Text  = "G1"
$t.Cell(5,9).Range.Text  = "GWY"
$t.Cell(5,10).Range.Text = "11.5"
$t.Cell(5,11).Range.Text = "0"
$t.Cell(5,12).Range.Text = "1"

# --- Total row (row 6) ---
$t.Cell(6,5).Range.Text = "8"
$t.Cell(6,6).Range.Text = "1"

# Merge the trailing empty cells of the 200m side (grid cols 7-12) the same
# way the 100m side already does: one gridSpan=4 "Total" label cell plus two
# value cells.
$totalCell = $t.Cell(6,7)
$totalCell.Merge($t.Cell(6,10))
$totalCell.Range.Text = "Total"
$totalCell.Range.ParagraphFormat.Alignment = 0

$t.Cell(6,11).Range.Text = "3"
$t.Cell(6,12).Range.Text = "6"
